$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mangled mobile phone number on row 2 (Q2 / Home_Phone)
$ws.Range("Q2").Value = "410-564-4639"

# Row 3: Mobile_Phone (P3) was numeric, now a formatted phone string
$ws.Range("P3").Value = "407-444-0909"

# Row 4: add missing last_name (E4) and fix Home_Phone (Q4)
$ws.Range("E4").Value = "Dumbface"
$ws.Range("Q4").Value = "407-721-7359"

# Row 4's auto-height shrinks slightly now that its content has changed
$ws.Rows.Item(4).RowHeight = 23.85

# Update the saved selection
$ws.Range("Q4").Select()
